$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "ful-path.csv" prediction scores in column B (rows 2-4) ---
$ws.Range("B2").Value = 1643.9452004067166
$ws.Range("B3").Value = 1535.0922864223048
$ws.Range("B4").Value = 1723.472007440319

# --- Re-stamp formatting on the header row and the row-label column so the
#     sheet picks up a fresh (but equivalent) cell style, mirroring the
#     restyle that happened when outputs-r202 was regenerated ---
$ws.Range("A1:C1").WrapText = $false
$ws.Range("A2:A4").WrapText = $false
